$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "id_type" masterdata sheet is being restructured: a new leading "id"
# column is introduced (A), the remaining columns are reordered/renamed
# (lang_code, code, name, descr, is_active), and the single sample data row
# is replaced with ten real id_type rows (English + French).
# ---------------------------------------------------------------------------

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "descr"
$ws.Range("F1").Value = "is_active"

# ---- Data rows (rows 2-11) ----
# columns: A=id, B=lang_code, C=code, D=name, E=descr, F=is_active

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "eng"
$ws.Range("C2").Value = "UIN"
$ws.Range("D2").Value = "Unique Identification Number"
$ws.Range("E2").Value = "National ID given to the applicant"
$ws.Range("F2").Value = $true

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "eng"
$ws.Range("C3").Value = "PRID"
$ws.Range("D3").Value = "Pre-registration ID"
$ws.Range("E3").Value = "ID assigned after Pre-registration"
$ws.Range("F3").Value = $true

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "eng"
$ws.Range("C4").Value = "RID"
$ws.Range("D4").Value = "Registration ID"
$ws.Range("E4").Value = "ID assigned after registration"
$ws.Range("F4").Value = $true

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "eng"
$ws.Range("C5").Value = "VID"
$ws.Range("D5").Value = "Virtual ID"
$ws.Range("E5").Value = "ID used in replacement of UIN"
$ws.Range("F5").Value = $true

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "eng"
$ws.Range("C6").Value = "Token ID"
$ws.Range("D6").Value = "Token ID"
$ws.Range("E6").Value = "ID used by a vendor for an applicant"
$ws.Range("F6").Value = $true

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "fra"
$ws.Range("C7").Value = "UIN"
$ws.Range("D7").Value = "Numéro didentification unique"
$ws.Range("E7").Value = "Carte didentité nationale fournie au demandeur"
$ws.Range("F7").Value = $true

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "fra"
$ws.Range("C8").Value = "PRID"
$ws.Range("D8").Value = "ID de pré-inscription"
$ws.Range("E8").Value = "ID attribué après la pré-inscription"
$ws.Range("F8").Value = $true

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "fra"
$ws.Range("C9").Value = "RID"
$ws.Range("D9").Value = "ID denregistrement"
$ws.Range("E9").Value = "ID attribué après lenregistrement"
$ws.Range("F9").Value = $true

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "fra"
$ws.Range("C10").Value = "VID"
$ws.Range("D10").Value = "ID virtuel"
$ws.Range("E10").Value = "Identifiant utilisé en remplacement de UIN"
$ws.Range("F10").Value = $true

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "fra"
$ws.Range("C11").Value = "Token ID"
$ws.Range("D11").Value = "ID de jeton"
$ws.Range("E11").Value = "ID utilisé par un fournisseur pour un demandeur"
$ws.Range("F11").Value = $true

# ---------------------------------------------------------------------------
# Formatting: the bold/bordered header style that originally lived on A1
# ("code") now belongs on the header row (B1:F1) and on the new id column
# (A2:A11). Propagate it via copy/paste-format *before* clearing A1, since
# A1 itself is unused by the new layout.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# A1 is no longer part of the table - clear its old content + style.
$ws.Range("A1").Clear()
